$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ================= Header text updates =================
$ws.Range("A8").Value = "Volume 31   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/15/2024  Through  7/21/2024"

# ================= Column width update =================
$ws.Columns("H").ColumnWidth = 6.168446

# ================= Row 15 =================
$ws.Range("N15").Value = -50

# ================= Row 16 =================
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 8.333333333333
$ws.Range("I16").Value = 134
$ws.Range("J16").Value = 75
$ws.Range("K16").Value = 78.666666666666
$ws.Range("L16").Value = 52.272727272727
$ws.Range("M16").Value = 21.818181818181
$ws.Range("N16").Value = -80.118694362017

# ================= Row 17 =================
$ws.Range("C17").Value = 6
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 14.814814814814
$ws.Range("I17").Value = 188
$ws.Range("J17").Value = 142
$ws.Range("K17").Value = 32.394366197183
$ws.Range("L17").Value = 86.138613861386
$ws.Range("M17").Value = 203.225806451613
$ws.Range("N17").Value = -15.695067264574

# ================= Row 18 =================
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 14
$ws.Range("I18").Value = 116
$ws.Range("J18").Value = 128
$ws.Range("K18").Value = -9.375
$ws.Range("L18").Value = 1.754385964912
$ws.Range("M18").Value = 81.25
$ws.Range("N18").Value = -74.107142857142

# ================= Row 19 =================
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 3.508771929824
$ws.Range("I19").Value = 404
$ws.Range("J19").Value = 383
$ws.Range("K19").Value = 5.483028720626
$ws.Range("L19").Value = 16.426512968299
$ws.Range("M19").Value = 59.683794466403
$ws.Range("N19").Value = -30.464716006884

# ================= Row 20 =================
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 17
$ws.Range("H20").Value = 183.333333333333
$ws.Range("I20").Value = 49
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 22.5
$ws.Range("L20").Value = 63.333333333333
$ws.Range("M20").Value = 81.481481481481
$ws.Range("N20").Value = -87.626262626262

# ================= Row 21 =================
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 136
$ws.Range("G21").Value = 116
$ws.Range("H21").Value = 17.241379310344
$ws.Range("I21").Value = 898
$ws.Range("J21").Value = 771
$ws.Range("K21").Value = 16.472114137483
$ws.Range("L21").Value = 31.094890510948
$ws.Range("M21").Value = 73.359073359073
$ws.Range("N21").Value = -61.558219178082

# ================= Row 22 =================
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 27
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = -10
$ws.Range("L22").Value = -27.027027027027
$ws.Range("M22").Value = -28.947368421052

# ================= Row 23 =================
$ws.Range("C23").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 25
$ws.Range("K23").Value = -3.846153846153
$ws.Range("L23").Value = 8.695652173913
$ws.Range("M23").Value = 92.307692307692

# ================= Row 24 =================
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 46
$ws.Range("E24").Value = 19.565217391304
$ws.Range("F24").Value = 196
$ws.Range("G24").Value = 170
$ws.Range("H24").Value = 15.294117647058
$ws.Range("I24").Value = 1215
$ws.Range("J24").Value = 1127
$ws.Range("K24").Value = 7.808340727595
$ws.Range("L24").Value = 23.727087576374
$ws.Range("M24").Value = 48.714810281517

# ================= Row 25 =================
$ws.Range("C25").Value = 48
$ws.Range("D25").Value = 49
$ws.Range("E25").Value = -2.040816326530
$ws.Range("F25").Value = 178
$ws.Range("G25").Value = 163
$ws.Range("H25").Value = 9.202453987730
$ws.Range("I25").Value = 1098
$ws.Range("J25").Value = 1017
$ws.Range("K25").Value = 7.964601769911
$ws.Range("L25").Value = 37.593984962406

# ================= Row 26 =================
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 33.333333333333
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 261
$ws.Range("J26").Value = 204
$ws.Range("K26").Value = 27.941176470588
$ws.Range("L26").Value = 30.5
$ws.Range("M26").Value = 32.48730964467

# ================= Row 27 =================
$ws.Range("D27").Value = 1
$ws.Range("F14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 80

# ================= Row 28 =================
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 33
$ws.Range("J28").Value = 29
$ws.Range("K28").Value = 13.793103448275
$ws.Range("L28").Value = -13.157894736842

# ================= Row 29 =================
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("H29").PasteSpecial(-4122)

# ================= Row 30 =================
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("H30").PasteSpecial(-4122)

# ================= Row 31 =================
$ws.Range("C31").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("F31").Value = 1
$ws.Range("F14").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 7
$ws.Range("K31").Value = -41.666666666666
$ws.Range("L31").Value = -22.222222222222

$excel.CutCopyMode = $false